$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was updated from
# 45192 to 45202 for every data row (rows 2 through 409).
$ws.Range("C2:C409").Value = 45202
